# Auto-generated edit script applying the diff to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.593.88"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.663.89"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.38"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.898.15"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.668.59"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.44%  "
$ws.Range("E14").Value = "  -2.90%  "
$ws.Range("E15").Value = "  -3.37%  "
$ws.Range("E16").Value = "  -3.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "246.71"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.579.64"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.46"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -7.07%  "
$ws.Range("E22").Value = "  -3.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.34"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.92%  "
$ws.Range("E24").Value = "  -4.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.24"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.58%  "
$ws.Range("E26").Value = "  -5.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.21"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  -2.03%  "
$ws.Range("E30").Value = "  +4.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0500"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.452.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("E34").Value = "  -5.10%  "
$ws.Range("E35").Value = "  -7.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.934"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.47%  "
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0170"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("E40").Value = "  -3.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.12"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.98%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.42"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -7.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.792"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.806.88"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.64%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.21"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.70"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.08"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0108"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.82"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.85%  "
